$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,7).Value = 319.073924
$ws.Cells.Item(2,8).Value = 957.221772
$ws.Cells.Item(2,9).Value = 0.6828333423212949
$ws.Cells.Item(2,10).Value = 0.6828333423212949
$ws.Cells.Item(2,13).Value = 40.81054266666667
$ws.Cells.Item(2,14).Value = 122.431628
$ws.Cells.Item(2,15).Value = 0.9943414173631485
$ws.Cells.Item(2,16).Value = 0.9943414173631485
$ws.Cells.Item(2,17).Value = 13021.57998922276
$ws.Cells.Item(2,18).Value = 117194.2199030048
$ws.Cells.Item(2,19).Value = 0.6789694734265723
$ws.Cells.Item(2,20).Value = 0.6789694734265723
$ws.Cells.Item(3,7).Value = 319.073924
$ws.Cells.Item(3,8).Value = 957.221772
$ws.Cells.Item(3,9).Value = 0.6828333423212949
$ws.Cells.Item(3,10).Value = 0.6828333423212949
$ws.Cells.Item(3,15).Value = 0.0002749163555820933
$ws.Cells.Item(3,16).Value = 0.0002749163555820933
$ws.Cells.Item(3,17).Value = 3.600217442466666
$ws.Cells.Item(3,18).Value = 32.4019569822
$ws.Cells.Item(3,19).Value = 0.0001877220539409103
$ws.Cells.Item(3,20).Value = 0.0001877220539409104
$ws.Cells.Item(4,7).Value = 319.073924
$ws.Cells.Item(4,8).Value = 957.221772
$ws.Cells.Item(4,9).Value = 0.6828333423212949
$ws.Cells.Item(4,10).Value = 0.6828333423212949
$ws.Cells.Item(4,15).Value = 0.002616675800765965
$ws.Cells.Item(4,16).Value = 0.002616675800765965
$ws.Cells.Item(4,17).Value = 34.26715678392933
$ws.Cells.Item(4,18).Value = 308.404411055364
$ws.Cells.Item(4,19).Value = 0.001786753482808274
$ws.Cells.Item(4,20).Value = 0.001786753482808275
$ws.Cells.Item(5,7).Value = 319.073924
$ws.Cells.Item(5,8).Value = 957.221772
$ws.Cells.Item(5,9).Value = 0.6828333423212949
$ws.Cells.Item(5,10).Value = 0.6828333423212949
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.113565
$ws.Cells.Item(5,14).Value = 0.340695
$ws.Cells.Item(5,15).Value = 0.002766990480503436
$ws.Cells.Item(5,16).Value = 0.002766990480503436
$ws.Cells.Item(5,17).Value = 36.23563017906
$ws.Cells.Item(5,18).Value = 326.12067161154
$ws.Cells.Item(5,19).Value = 0.001889393357973367
$ws.Cells.Item(5,20).Value = 0.001889393357973367
$ws.Cells.Item(6,7).Value = 140.4344916666667
$ws.Cells.Item(6,8).Value = 421.303475
$ws.Cells.Item(6,9).Value = 0.3005364779415257
$ws.Cells.Item(6,10).Value = 0.3005364779415257
$ws.Cells.Item(6,13).Value = 40.81054266666667
$ws.Cells.Item(6,14).Value = 122.431628
$ws.Cells.Item(6,15).Value = 0.9943414173631485
$ws.Cells.Item(6,16).Value = 0.9943414173631485
$ws.Cells.Item(6,17).Value = 5731.207814034145
$ws.Cells.Item(6,18).Value = 51580.8703263073
$ws.Cells.Item(6,19).Value = 0.2988358674457052
$ws.Cells.Item(6,20).Value = 0.2988358674457052
$ws.Cells.Item(7,7).Value = 140.4344916666667
$ws.Cells.Item(7,8).Value = 421.303475
$ws.Cells.Item(7,9).Value = 0.3005364779415257
$ws.Cells.Item(7,10).Value = 0.3005364779415257
$ws.Cells.Item(7,15).Value = 0.0002749163555820933
$ws.Cells.Item(7,16).Value = 0.0002749163555820933
$ws.Cells.Item(7,19).Value = 0.0000826223932351624
$ws.Cells.Item(7,20).Value = 0.00008262239323516242
$ws.Cells.Item(8,7).Value = 140.4344916666667
$ws.Cells.Item(8,8).Value = 421.303475
$ws.Cells.Item(8,9).Value = 0.3005364779415257
$ws.Cells.Item(8,10).Value = 0.3005364779415257
$ws.Cells.Item(8,15).Value = 0.002616675800765965
$ws.Cells.Item(8,16).Value = 0.002616675800765965
$ws.Cells.Item(8,19).Value = 0.0007864065290770243
$ws.Cells.Item(8,20).Value = 0.0007864065290770244
$ws.Cells.Item(9,7).Value = 140.4344916666667
$ws.Cells.Item(9,8).Value = 421.303475
$ws.Cells.Item(9,9).Value = 0.3005364779415257
$ws.Cells.Item(9,10).Value = 0.3005364779415257
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.113565
$ws.Cells.Item(9,14).Value = 0.340695
$ws.Cells.Item(9,15).Value = 0.002766990480503436
$ws.Cells.Item(9,16).Value = 0.002766990480503436
$ws.Cells.Item(9,17).Value = 15.948443046125
$ws.Cells.Item(9,18).Value = 143.535987415125
$ws.Cells.Item(9,19).Value = 0.0008315815735082323
$ws.Cells.Item(9,20).Value = 0.0008315815735082323
$ws.Cells.Item(10,7).Value = 7.479044333333333
$ws.Cells.Item(10,8).Value = 22.437133
$ws.Cells.Item(10,9).Value = 0.01600550986892662
$ws.Cells.Item(10,10).Value = 0.01600550986892662
$ws.Cells.Item(10,13).Value = 40.81054266666667
$ws.Cells.Item(10,14).Value = 122.431628
$ws.Cells.Item(10,15).Value = 0.9943414173631485
$ws.Cells.Item(10,16).Value = 0.9943414173631485
$ws.Cells.Item(10,17).Value = 305.2238578713916
$ws.Cells.Item(10,18).Value = 2747.014720842524
$ws.Cells.Item(10,19).Value = 0.01591494136868835
$ws.Cells.Item(10,20).Value = 0.01591494136868835
$ws.Cells.Item(11,7).Value = 7.479044333333333
$ws.Cells.Item(11,8).Value = 22.437133
$ws.Cells.Item(11,9).Value = 0.01600550986892662
$ws.Cells.Item(11,10).Value = 0.01600550986892662
$ws.Cells.Item(11,15).Value = 0.0002749163555820933
$ws.Cells.Item(11,16).Value = 0.0002749163555820933
$ws.Cells.Item(11,17).Value = 0.08438855022777778
$ws.Cells.Item(11,18).Value = 0.7594969520499999
$ws.Cells.Item(11,19).Value = 0.000004400176442398533
$ws.Cells.Item(11,20).Value = 0.000004400176442398534
$ws.Cells.Item(12,7).Value = 7.479044333333333
$ws.Cells.Item(12,8).Value = 22.437133
$ws.Cells.Item(12,9).Value = 0.01600550986892662
$ws.Cells.Item(12,10).Value = 0.01600550986892662
$ws.Cells.Item(12,15).Value = 0.002616675800765965
$ws.Cells.Item(12,16).Value = 0.002616675800765965
$ws.Cells.Item(12,17).Value = 0.8032169522078889
$ws.Cells.Item(12,18).Value = 7.228952569871
$ws.Cells.Item(12,19).Value = 0.0000418812303529411
$ws.Cells.Item(12,20).Value = 0.00004188123035294111
$ws.Cells.Item(13,7).Value = 7.479044333333333
$ws.Cells.Item(13,8).Value = 22.437133
$ws.Cells.Item(13,9).Value = 0.01600550986892662
$ws.Cells.Item(13,10).Value = 0.01600550986892662
$ws.Cells.Item(13,11).Value = 2
$ws.Cells.Item(13,12).Value = 0.6666666666666666
$ws.Cells.Item(13,13).Value = 0.113565
$ws.Cells.Item(13,14).Value = 0.340695
$ws.Cells.Item(13,15).Value = 0.002766990480503436
$ws.Cells.Item(13,16).Value = 0.002766990480503436
$ws.Cells.Item(13,17).Value = 0.8493576697150002
$ws.Cells.Item(13,18).Value = 7.644219027435001
$ws.Cells.Item(13,19).Value = 0.00004428709344292374
$ws.Cells.Item(13,20).Value = 0.00004428709344292374
$ws.Cells.Item(14,5).Value = 1
$ws.Cells.Item(14,6).Value = 0.3333333333333333
$ws.Cells.Item(14,7).Value = 0.2918953333333333
$ws.Cells.Item(14,8).Value = 0.875686
$ws.Cells.Item(14,9).Value = 0.0006246698682528143
$ws.Cells.Item(14,10).Value = 0.0006246698682528143
$ws.Cells.Item(14,13).Value = 40.81054266666667
$ws.Cells.Item(14,14).Value = 122.431628
$ws.Cells.Item(14,15).Value = 0.9943414173631485
$ws.Cells.Item(14,16).Value = 0.9943414173631485
$ws.Cells.Item(14,17).Value = 11.91240695520089
$ws.Cells.Item(14,18).Value = 107.211662596808
$ws.Cells.Item(14,19).Value = 0.0006211351221825546
$ws.Cells.Item(14,20).Value = 0.0006211351221825546
$ws.Cells.Item(15,5).Value = 1
$ws.Cells.Item(15,6).Value = 0.3333333333333333
$ws.Cells.Item(15,7).Value = 0.2918953333333333
$ws.Cells.Item(15,8).Value = 0.875686
$ws.Cells.Item(15,9).Value = 0.0006246698682528143
$ws.Cells.Item(15,10).Value = 0.0006246698682528143
$ws.Cells.Item(15,15).Value = 0.0002749163555820933
$ws.Cells.Item(15,16).Value = 0.0002749163555820933
$ws.Cells.Item(15,17).Value = 0.003293552344444444
$ws.Cells.Item(15,18).Value = 0.0296419711
$ws.Cells.Item(15,19).Value = 0.0000001717319636220101
$ws.Cells.Item(15,20).Value = 0.0000001717319636220101
$ws.Cells.Item(16,5).Value = 1
$ws.Cells.Item(16,6).Value = 0.3333333333333333
$ws.Cells.Item(16,7).Value = 0.2918953333333333
$ws.Cells.Item(16,8).Value = 0.875686
$ws.Cells.Item(16,9).Value = 0.0006246698682528143
$ws.Cells.Item(16,10).Value = 0.0006246698682528143
$ws.Cells.Item(16,15).Value = 0.002616675800765965
$ws.Cells.Item(16,16).Value = 0.002616675800765965
$ws.Cells.Item(16,17).Value = 0.03134829392022222
$ws.Cells.Item(16,18).Value = 0.282134645282
$ws.Cells.Item(16,19).Value = 0.000001634558527724802
$ws.Cells.Item(16,20).Value = 0.000001634558527724803
$ws.Cells.Item(17,5).Value = 1
$ws.Cells.Item(17,6).Value = 0.3333333333333333
$ws.Cells.Item(17,7).Value = 0.2918953333333333
$ws.Cells.Item(17,8).Value = 0.875686
$ws.Cells.Item(17,9).Value = 0.0006246698682528143
$ws.Cells.Item(17,10).Value = 0.0006246698682528143
$ws.Cells.Item(17,11).Value = 2
$ws.Cells.Item(17,12).Value = 0.6666666666666666
$ws.Cells.Item(17,13).Value = 0.113565
$ws.Cells.Item(17,14).Value = 0.340695
$ws.Cells.Item(17,15).Value = 0.002766990480503436
$ws.Cells.Item(17,16).Value = 0.002766990480503436
$ws.Cells.Item(17,17).Value = 0.03314909353000001
$ws.Cells.Item(17,18).Value = 0.29834184177
$ws.Cells.Item(17,19).Value = 0.000001728455578912872
$ws.Cells.Item(17,20).Value = 0.000001728455578912872
